$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
try {
    $wb2 = $excel.Workbooks.Add()
    $ws2 = $wb2.ActiveSheet
    Write-Host "New wb font color:" $ws2.Range("A1").Font.Color
    Write-Host "New wb font themecolor:" $ws2.Range("A1").Font.ThemeColor
    $ws.Range("B52").Value = 572
    $ws2.Range("A1").Copy()
    $ws.Range("B52").PasteSpecial(-4122)
} catch {
    Write-Host "ERR: $_"
}
